$wb = $excel.ActiveWorkbook

# --- Shared-string play-log appends (space-separated running totals) ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = $ws.Range("B2").Value() + " 3 -1 1 2 1 4 2 4 1 12 8 2 0 26 4 1 0 2 4 5 1 6 0 6 3 -1 6"
$ws.Range("C2").Value = $ws.Range("C2").Value() + " 8 1 2 3 0 10 -2 6 5 2 8 23 1 13 8"
$ws.Range("B3").Value = $ws.Range("B3").Value() + " 17 6 12 16 8 15 7 2 9 1 11 6 18 29 7 10 14 12 8 1 6 5 3"
$ws.Range("C3").Value = $ws.Range("C3").Value() + " 11 33 6 9 11 4 7 10 3 18 4 10 8 21 9 12"

$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 172
$ws.Range("F2").Value = 61
$ws.Range("G2").Value = 56
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 26
$ws.Range("L2").Value = 304
$ws.Range("M2").Value = 185
$ws.Range("O2").Value = 31
$ws.Range("P2").Value = 18
$ws.Range("Q2").Value = 548
$ws.Range("B3").Value = 26
$ws.Range("C3").Value = 181
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 35
$ws.Range("F3").Value = 96
$ws.Range("G3").Value = 43
$ws.Range("H3").Value = 32
$ws.Range("I3").Value = 56
$ws.Range("J3").Value = 57
$ws.Range("N3").Value = 18

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 229
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 69
$ws.Range("G2").Value = 72
$ws.Range("J2").Value = 37
$ws.Range("L2").Value = 294
$ws.Range("M2").Value = 204
$ws.Range("O2").Value = 20
$ws.Range("P2").Value = 14
$ws.Range("Q2").Value = 557
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 161
$ws.Range("E3").Value = 34
$ws.Range("F3").Value = 108
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 32
$ws.Range("I3").Value = 59
$ws.Range("J3").Value = 68
$ws.Range("N3").Value = 20

$ws = $wb.Worksheets.Item("ST")
$ws.Range("D3").Value = $ws.Range("D3").Value() + " 45 50 57"
$ws.Range("B4").Value = $ws.Range("B4").Value() + " 66 65 61"
$ws.Range("D4").Value = $ws.Range("D4").Value() + " 0 0 7"
$ws.Range("B5").Value = $ws.Range("B5").Value() + " 40 18 41"
$ws.Range("D5").Value = $ws.Range("D5").Value() + " 0 0"
$ws.Range("B2").Value = 59
$ws.Range("D2").Value = 70
$ws.Range("F2").Value = 22
$ws.Range("J2").Value = 52
$ws.Range("K2").Value = 46
$ws.Range("B3").Value = 22

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value = 7
$ws.Range("E2").Value = 6

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 26
$ws.Range("B3").Value = 19
$ws.Range("B4").Value = 4
$ws.Range("D4").Value = 4

